$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (StatQuery), shifting old B->C and C->D
$ws.Columns.Item(2).Insert()

# Match column width of the new column B to column A's width (75.81640625 chars)
$ws.Columns.Item(2).ColumnWidth = 75

# Header row
$ws.Range("B1").Value = "StatQuery"

# Data row - new stat query text, with the same wrap-text style used by A2
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE t.clinical_trial_id IN ['NCT02465060'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B2").WrapText = $true

# Update the selected range to the whole of the new column B, matching author's view state
$null = $ws.Range("B:B").Select()
